$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date (D2:D6) from 2021-12-06 to 2021-12-07
$ws.Range("D2:D6").Value = 44537

# Update time_last_updated (E2:E6) epoch seconds
$ws.Range("E2:E6").Value = 1638835201

# Update rate values that changed (F4, F5)
$ws.Range("F4").Value = 107.1
$ws.Range("F5").Value = 490.9
